# Apply updated pricing/profit figures (scheduled runner refresh)
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H11").Value = 59.714287
$ws.Range("I11").Value = 59.714287
$ws.Range("K11").Value = 59.714287
$ws.Range("M11").Value = 80.285713

$ws.Range("H28").Value = 1209.24
$ws.Range("I28").Value = 1076.9474
$ws.Range("J28").Value = 1628.1666
$ws.Range("K28").Value = 1076.9474
$ws.Range("L28").Value = 1628.1666
$ws.Range("M28").Value = -591.9474
$ws.Range("N28").Value = -2598.1666

$ws.Range("H107").Value = 653.84
$ws.Range("I107").Value = 658
$ws.Range("J107").Value = 623.3333
$ws.Range("K107").Value = 658
$ws.Range("L107").Value = 623.3333
$ws.Range("M107").Value = 1262
$ws.Range("N107").Value = -4463.3333

$ws.Range("H111").Value = 2321.8
$ws.Range("I111").Value = 1707.091
$ws.Range("J111").Value = 4012.25
$ws.Range("K111").Value = 5121.272999999999
$ws.Range("L111").Value = 12036.75
$ws.Range("M111").Value = -2054.272999999999
$ws.Range("N111").Value = -18170.75

$ws.Range("H113").Value = 13192.357
$ws.Range("I113").Value = 6296.5
$ws.Range("K113").Value = 6296.5
$ws.Range("M113").Value = -3042.5

$ws.Range("H125").Value = 1523.75
$ws.Range("I125").Value = 1470
$ws.Range("J125").Value = 1900
$ws.Range("K125").Value = 13230
$ws.Range("L125").Value = 17100
$ws.Range("M125").Value = -10770
$ws.Range("N125").Value = -22020

$ws.Range("H137").Value = 4409.08
$ws.Range("J137").Value = 5117.625
$ws.Range("L137").Value = 15352.875
$ws.Range("N137").Value = -20452.875

$ws.Range("H138").Value = 2073.0857
$ws.Range("I138").Value = 1410.9584
$ws.Range("J138").Value = 3517.7273
$ws.Range("K138").Value = 4232.8752
$ws.Range("L138").Value = 10553.1819
$ws.Range("M138").Value = 907.1247999999996
$ws.Range("N138").Value = -20833.1819

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 36383.3
$ws.Range("I2").Value = 56337
$ws.Range("J2").Value = 6452.75
$ws.Range("K2").Value = 56337
$ws.Range("L2").Value = 6452.75
$ws.Range("M2").Value = -56224
$ws.Range("N2").Value = -6678.75

$ws.Range("H32").Value = 2394
$ws.Range("I32").Value = 2515.6052
$ws.Range("J32").Value = 1623.8334
$ws.Range("K32").Value = 2515.6052
$ws.Range("L32").Value = 1623.8334
$ws.Range("M32").Value = -2228.6052
$ws.Range("N32").Value = -2197.8334

$ws.Range("H61").Value = 11419.1875
$ws.Range("I61").Value = 3437.9092
$ws.Range("K61").Value = 3437.9092
$ws.Range("M61").Value = -3225.9092

$ws.Range("H88").Value = 1497.9412
$ws.Range("J88").Value = 1399.4286
$ws.Range("L88").Value = 1399.4286
$ws.Range("N88").Value = -2211.4286

$ws.Range("H91").Value = 1497.9412
$ws.Range("J91").Value = 1399.4286
$ws.Range("L91").Value = 1399.4286
$ws.Range("N91").Value = -4207.4286

$ws.Range("H102").Value = 7296.3184
$ws.Range("I102").Value = 2922.1052
$ws.Range("J102").Value = 34999.668
$ws.Range("K102").Value = 2922.1052
$ws.Range("L102").Value = 34999.668
$ws.Range("M102").Value = -1300.1052
$ws.Range("N102").Value = -38243.668

$ws.Range("H110").Value = 6023.7
$ws.Range("I110").Value = 5974.857
$ws.Range("K110").Value = 5974.857
$ws.Range("M110").Value = -3929.857

$ws.Range("H116").Value = 36383.3
$ws.Range("I116").Value = 56337
$ws.Range("J116").Value = 6452.75
$ws.Range("K116").Value = 56337
$ws.Range("L116").Value = 6452.75
$ws.Range("M116").Value = -54043
$ws.Range("N116").Value = -11040.75

$ws.Range("H122").Value = 3581
$ws.Range("I122").Value = 2948.6667
$ws.Range("K122").Value = 8846.000100000001
$ws.Range("M122").Value = -6396.000100000001

$ws.Range("H136").Value = 11419.1875
$ws.Range("I136").Value = 3437.9092
$ws.Range("K136").Value = 10313.7276
$ws.Range("M136").Value = -7763.7276

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 36383.3
$ws.Range("I3").Value = 56337
$ws.Range("J3").Value = 6452.75
$ws.Range("K3").Value = 56337
$ws.Range("L3").Value = 6452.75
$ws.Range("M3").Value = -56223
$ws.Range("N3").Value = -6680.75

$ws.Range("H105").Value = 2002
$ws.Range("I105").Value = 1083.9231
$ws.Range("J105").Value = 3707
$ws.Range("K105").Value = 1083.9231
$ws.Range("L105").Value = 3707
$ws.Range("M105").Value = 663.0769
$ws.Range("N105").Value = -7201

$ws.Range("H107").Value = 12456.143
$ws.Range("I107").Value = 11496.1
$ws.Range("K107").Value = 11496.1
$ws.Range("M107").Value = -9576.1

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H99").Value = 2561.8667
$ws.Range("I99").Value = 2297.9
$ws.Range("J99").Value = 3089.8
$ws.Range("K99").Value = 2297.9
$ws.Range("L99").Value = 3089.8
$ws.Range("M99").Value = -799.9000000000001
$ws.Range("N99").Value = -6085.8

$ws.Range("H107").Value = 5501.95
$ws.Range("I107").Value = 228.72728
$ws.Range("J107").Value = 11947
$ws.Range("K107").Value = 228.72728
$ws.Range("L107").Value = 11947
$ws.Range("M107").Value = 1691.27272
$ws.Range("N107").Value = -15787

$ws.Range("H122").Value = 3574
$ws.Range("I122").Value = 2532.6667
$ws.Range("J122").Value = 4912.857
$ws.Range("K122").Value = 7598.000100000001
$ws.Range("L122").Value = 14738.571
$ws.Range("M122").Value = -5148.000100000001
$ws.Range("N122").Value = -19638.571

$ws.Range("H126").Value = 2561.8667
$ws.Range("I126").Value = 2297.9
$ws.Range("J126").Value = 3089.8
$ws.Range("K126").Value = 6893.700000000001
$ws.Range("L126").Value = 9269.400000000001
$ws.Range("M126").Value = -4423.700000000001
$ws.Range("N126").Value = -14209.4

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 1541359.1
$ws.Range("J131").Value = 5760
$ws.Range("L131").Value = 17280
$ws.Range("N131").Value = -27360

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H38").Value = 22997
$ws.Range("J38").Value = 22997
$ws.Range("L38").Value = 22997
$ws.Range("N38").Value = -23923

$ws.Range("H70").Value = 80382.87
$ws.Range("I70").Value = 107263.27
$ws.Range("J70").Value = 6461.75
$ws.Range("K70").Value = 107263.27
$ws.Range("L70").Value = 6461.75
$ws.Range("M70").Value = -106993.27
$ws.Range("N70").Value = -7001.75

$ws.Range("H73").Value = 80382.87
$ws.Range("I73").Value = 107263.27
$ws.Range("J73").Value = 6461.75
$ws.Range("K73").Value = 107263.27
$ws.Range("L73").Value = 6461.75
$ws.Range("M73").Value = -106327.27
$ws.Range("N73").Value = -8333.75

$ws.Range("H102").Value = 4772.294
$ws.Range("I102").Value = 3945.5625
$ws.Range("K102").Value = 3945.5625
$ws.Range("M102").Value = -2323.5625

$ws.Range("H107").Value = 341.83334
$ws.Range("I107").Value = 341.83334
$ws.Range("J107").Value = 0
$ws.Range("K107").Value = 341.83334
$ws.Range("L107").Value = 0
$ws.Range("M107").Value = 1578.16666
$ws.Range("N107").ClearContents()

$ws.Range("H132").Value = 8090.0967
$ws.Range("I132").Value = 7377.9565
$ws.Range("J132").Value = 10137.5
$ws.Range("K132").Value = 22133.8695
$ws.Range("L132").Value = 30412.5
$ws.Range("M132").Value = -19603.8695
$ws.Range("N132").Value = -35472.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H55").Value = 446.32
$ws.Range("I55").Value = 444.46667
$ws.Range("K55").Value = 444.46667
$ws.Range("M55").Value = -271.46667

$ws.Range("H100").Value = 124844.445
$ws.Range("I100").Value = 160301.42
$ws.Range("K100").Value = 160301.42
$ws.Range("M100").Value = -159760.42

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 753.6667
$ws.Range("I107").Value = 752.6429000000001
$ws.Range("K107").Value = 2257.9287
$ws.Range("M107").Value = -337.9287000000004

$ws.Range("H126").Value = 3593
$ws.Range("I126").Value = 3593
$ws.Range("K126").Value = 10779
$ws.Range("M126").Value = -8309

$ws.Range("H132").Value = 4743.6665
$ws.Range("I132").Value = 3783.125
$ws.Range("J132").Value = 6664.75
$ws.Range("K132").Value = 11349.375
$ws.Range("L132").Value = 19994.25
$ws.Range("M132").Value = -8819.375
$ws.Range("N132").Value = -25054.25
